$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new timesheet entry row (2014-12-02) right after the last data row,
# copying row 22's formatting (date / time / general styles) so the new
# row matches the existing table's look.
$ws.Range("A22:D22").Copy()
$ws.Range("A23:D23").PasteSpecial(-4122)

$ws.Range("A23").Value = 41975
$ws.Range("B23").Value = 0.0833333333333333
$ws.Range("C23").Value = 0.229166666666667
$ws.Range("D23").Formula = "=ROUND(ABS(C23-B23) * 24, 1)"

# Leave the selection where the user last clicked while editing.
$ws.Range("C23").Select()
